# Fix for issue: XLSReader tends to convert long number to scientific notation
#
# Adds a new test row (row 11) to Sheet1 demonstrating a long integer value
# that should be preserved verbatim (not converted to scientific notation):
#   A11 = 10
#   B11 = 123456789012345   (the long number in question)
#   C11 = "gagag"           (new shared string)
# and leaves the selection on the newly added cell B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 123456789012345
$ws.Range("C11").Value = "gagag"

# Matches the author's final selection state (B11) recorded in the sheet XML.
$ws.Range("B11").Select() | Out-Null
